$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.601.14"
$ws.Range("E2").Value = "  +4.94%  "
$ws.Range("D3").Value = "3.624.52"
$ws.Range("E3").Value = "  +18.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.48"
$ws.Range("E5").Value = "  +3.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.62"
$ws.Range("E6").Value = "  +8.63%  "
$ws.Range("D7").Value = "3.625.83"
$ws.Range("E7").Value = "  +18.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  +4.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  +7.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.51"
$ws.Range("E11").Value = "  +3.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.495"
$ws.Range("E12").Value = "  +5.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.97"
$ws.Range("E13").Value = "  +9.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000252"
$ws.Range("E14").Value = "  +5.94%  "
$ws.Range("D15").Value = "4.226.10"
$ws.Range("E15").Value = "  +17.95%  "
$ws.Range("D16").Value = "3.623.90"
$ws.Range("E16").Value = "  +18.03%  "
$ws.Range("D17").Value = "69.748.61"
$ws.Range("E17").Value = "  +5.32%  "
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.53"
$ws.Range("E19").Value = "  +8.39%  "
$ws.Range("E20").Value = "  +3.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "506.03"
$ws.Range("E21").Value = "  +4.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.31"
$ws.Range("E22").Value = "  +21.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.744"
$ws.Range("E23").Value = "  +8.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.09"
$ws.Range("E24").Value = "  +5.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.46"
$ws.Range("E25").Value = "  +6.71%  "
$ws.Range("E26").Value = "  +8.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.83"
$ws.Range("E27").Value = "  +6.89%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +12.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.10"
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.81"
$ws.Range("E31").Value = "  +18.63%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.75"
$ws.Range("E32").Value = "  +5.85%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0000109"
$ws.Range("E33").Value = "  +20.53%  "
$ws.Range("E34").Value = "  +5.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.13"
$ws.Range("E36").Value = "  +10.94%  "
$ws.Range("E37").Value = "  +8.70%  "
$ws.Range("E38").Value = "  +11.41%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.11"
$ws.Range("E39").Value = "  +7.90%  "
$ws.Range("B40").Value = "Arweave"
$ws.Range("C40").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "46.61"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.62"
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("E42").Value = "  +4.70%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.125.34"
$ws.Range("E43").Value = "  +13.09%  "
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.79"
$ws.Range("E44").Value = "  +7.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.80"
$ws.Range("E45").Value = "  +10.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "404.85"
$ws.Range("E46").Value = "  +11.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0366"
$ws.Range("E47").Value = "  +6.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.75"
$ws.Range("E48").Value = "  +14.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "136.21"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.44"
$ws.Range("E51").Value = "  +13.99%  "
